# Apply the "VENDA 08 (22/09)" batch update to the Cliente sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the access period (column E, "PRAZO ACESSO - DIAS") to 365 for
#     the active rows 5-12 ---
$ws.Range("E5").Value = 365
$ws.Range("E6").Value = 365
$ws.Range("E7").Value = 365
$ws.Range("E8").Value = 365
$ws.Range("E9").Value = 365
$ws.Range("E10").Value = 365
$ws.Range("E11").Value = 365
$ws.Range("E12").Value = 365

# --- Fill in the EMAIL column ("-") for rows 9-12, matching the same
#     formatting already used on rows 6-8 (column F, style carries the
#     bordered/filled "data row" look) ---
$ws.Range("F6").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = "-"

$ws.Range("F6").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value = "-"

$ws.Range("F6").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value = "-"

$ws.Range("F6").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = "-"

# --- Row 12: move the start date out to 44830 and mark it CONFIRMADO as
#     "VENDA 08 (22/09)" (new batch, replacing the old NILTON BISPO / TESTE
#     EDSON placeholder rows) ---
$ws.Range("D12").Value = 44830
$ws.Range("G9").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G12").Value = "VENDA 08 (22/09)"

# --- Row 13 (previously NILTON BISPO) is cleared out entirely, keeping the
#     existing cell formatting ---
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()

# --- Row 16 (previously TESTE EDSON) also loses its name/date/prazo data ---
$ws.Range("B16").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
